$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 39; $r++) {
    $ws.Cells.Item($r, 12).Value2 = $ws.Cells.Item($r, 11).Value2
}

$ws.Range("C53").Value2 = 45
